$wb = $excel.ActiveWorkbook

# --- Sheet 1: Coefficients ---
$wsCoef = $wb.Worksheets.Item("Coefficients")

# Replace "Cruise" labels with "Month" labels (cruise -> month update)
$wsCoef.Range("A5").Value = "MonthOctober"
$wsCoef.Range("A7").Value = "DRM:MonthOctober"
$wsCoef.Range("A8").Value = "Depth:MonthOctober"

# --- Sheet 2: Fullmodel_statistics ---
$wsStats = $wb.Worksheets.Item("Fullmodel_statistics")

# Minor floating point precision refresh on recomputed model statistics
$wsStats.Range("A2").Value = 0.4053986377333108
$wsStats.Range("B2").Value = 0.2681829387486901
$wsStats.Range("C2").Value = 2.954462504897116
